$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J, matching the style of the existing
# header row (e.g. H1 "IP": bold font, thin box border, centered/top aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy H1's cell format (bold font, border, center/top alignment) onto the
# new header cells, the same way a user would drag-fill/paste-format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill data rows 2-34: column I is always 1, column J mirrors column H.
$lastRow = 34
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
